# Shift the "Time Reduction" (98%) and "Document Parsing" (60s) metric
# text boxes on Slide 3 to the left, matching the updated layout.
#
# EMU -> point conversion used by the PowerPoint object model: 1 pt = 12700 EMU
#   1371600 EMU -> 108.0 pt   (old X for the first metric pair)
#   1097280 EMU ->  86.4 pt   (new X for the first metric pair)
#   4114800 EMU -> 324.0 pt   (old X for the second metric pair)
#   3931920 EMU -> 309.6 pt   (new X for the second metric pair)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# "Text 20" -> "98%" value box
$s.Shapes.Item("Text 20").Left = 86.4
# "Text 21" -> "Time Reduction" caption box
$s.Shapes.Item("Text 21").Left = 86.4
# "Text 22" -> "60s" value box
$s.Shapes.Item("Text 22").Left = 309.6
# "Text 23" -> "Document Parsing" caption box
$s.Shapes.Item("Text 23").Left = 309.6
